$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host $ws.Name
$ws.Rows("859:864").Delete()
Write-Host $ws.UsedRange.Address()
Write-Host $ws.Cells.Item(858,1).Value2
Write-Host $ws.Cells.Item(859,1).Value2
